$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values
$ws.Range("E11").Value = 27580
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

# Delete row 17 entirely (shifts rows below up by one)
$ws.Rows("17:17").Delete()
